# Generate Report for Handoff
#
# The localization status report just finished generating the handoff
# package: flip the "In Translation" status to "Ready for handoff" for
# every sheet that tracks it, bump the handoff timestamps to the new
# generation time, and widen the now-longer status/date columns so the
# new text isn't truncated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$overview.Range("F2").Value = "Ready for handoff"   # de-de status column
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status column
$dede.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Latest Handoff / Xliff-generate timestamps ---
$overview.Range("G2").Value = "2016-09-01 04:43:40" # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-09-01 04:43:40" # Latest Handoff Datetime (de-de)
$zhcn.Range("H2").Value     = "2016-09-01 04:43:35" # Latest Handoff Datetime (zh-cn)

# --- Column widths: widen the status/date columns for the longer text ---
$overview.Columns.Item(5).ColumnWidth = 16.38
$overview.Columns.Item(6).ColumnWidth = 16.38
$zhcn.Columns.Item(3).ColumnWidth     = 16.38
$dede.Columns.Item(3).ColumnWidth     = 16.38
